$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 22 (shifts old rows 22:133 down to 23:134), scoped
# to columns A:K only so unrelated columns aren't touched.
[void]$ws.Range("A22:K22").Insert(-4121)

# The freshly inserted row picks up brand-new default formatting; restore the
# normal "table body row" look by copying formats from the row directly
# below (the old row 22, now shifted to row 23), which still carries the
# correct original styles.
[void]$ws.Range("A23:K23").Copy()
[void]$ws.Range("A22:K22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-assert the calculated-column helper formula in the new row.
$ws.Range("G22").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# The row that used to be the table's last row (133) now lives on row 134,
# which briefly sat outside Table1's bounds during the shift; its
# calculated-column formula got rewritten as an unqualified structured
# reference ([@EARNED]) that errors (#VALUE!) outside a table context.
# Restore the proper formula text.
$ws.Range("G134").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# New leave-card entry: a 2-day "SL(2-0-0)" sick-leave debit dated 3/2,3/2023.
$ws.Range("B22").Value = "SL(2-0-0)"
$ws.Range("H22").Value = 2

# The K column on a "remarks with date" row (like K21) uses a dedicated
# date-flavoured style, so pull its format before writing the remark text.
[void]$ws.Range("K21").Copy()
[void]$ws.Range("K22").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K22").Value = "3/2,3/2023"

# Grow Table1 so it covers the newly inserted row.
$lo = $ws.ListObjects.Item("Table1")
[void]$lo.Resize($ws.Range("A8:K134"))

# Match the author's final cell selection.
[void]$ws.Range("B23").Select()

[void]$wb.Application.Calculate()
